$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use existing row 129 (a fully populated row) as the formatting template
# for the new rows so that styles (bold/border id column, date format
# column, etc.) match the rest of the sheet exactly.
$ws.Range("A129:AC129").Copy()
$ws.Range("A131:AC134").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Row 131 ----
$ws.Range("A131").Value = 129
$ws.Range("B131").Value = 7011624
$ws.Range("C131").Value = "Azerbaijan Premier League"
$ws.Range("D131").Value = "Azerbaijan Premier League"
$ws.Range("E131").Value = 45366.5
$ws.Range("F131").Value = "Neftchi Baku"
$ws.Range("G131").Value = "Sabail FC"
$ws.Range("H131").Value = 3
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = "H"
$ws.Range("K131").Value = 2.1
$ws.Range("L131").Value = 3.1
$ws.Range("M131").Value = 3.2
$ws.Range("N131").Value = 1.85
$ws.Range("O131").Value = 3.4
$ws.Range("P131").Value = 3.5
$ws.Range("Q131").Value = -0.5
$ws.Range("R131").Value = 1.95
$ws.Range("S131").Value = 1.85
$ws.Range("T131").Value = 2.5
$ws.Range("U131").Value = 1.775
$ws.Range("V131").Value = 1.925
$ws.Range("W131").Value = 0.8500000000000001
$ws.Range("X131").Value = -1
$ws.Range("Y131").Value = -1
$ws.Range("Z131").Value = 0.95
$ws.Range("AA131").Value = -1
$ws.Range("AB131").Value = 0.7749999999999999
$ws.Range("AC131").Value = -1

# ---- Row 132 ----
$ws.Range("A132").Value = 130
$ws.Range("B132").Value = 7011625
$ws.Range("C132").Value = "Azerbaijan Premier League"
$ws.Range("D132").Value = "Azerbaijan Premier League"
$ws.Range("E132").Value = 45367.375
$ws.Range("F132").Value = "Zira IK"
$ws.Range("G132").Value = "FK Kapaz"
$ws.Range("H132").Value = 2
$ws.Range("I132").Value = 1
$ws.Range("J132").Value = "H"
$ws.Range("K132").Value = 1.85
$ws.Range("L132").Value = 3.1
$ws.Range("M132").Value = 4
$ws.Range("N132").Value = 1.727
$ws.Range("O132").Value = 3.2
$ws.Range("P132").Value = 4.75
$ws.Range("Q132").Value = -0.75
$ws.Range("R132").Value = 1.975
$ws.Range("S132").Value = 1.825
$ws.Range("T132").Value = 2.25
$ws.Range("U132").Value = 1.975
$ws.Range("V132").Value = 1.725
$ws.Range("W132").Value = 0.7270000000000001
$ws.Range("X132").Value = -1
$ws.Range("Y132").Value = -1
$ws.Range("Z132").Value = 0.4875
$ws.Range("AA132").Value = -0.5
$ws.Range("AB132").Value = 0.9750000000000001
$ws.Range("AC132").Value = -1

# ---- Row 133 ----
$ws.Range("A133").Value = 131
$ws.Range("B133").Value = 7011622
$ws.Range("C133").Value = "Azerbaijan Premier League"
$ws.Range("D133").Value = "Azerbaijan Premier League"
$ws.Range("E133").Value = 45367.47916666666
$ws.Range("F133").Value = "Araz FK"
$ws.Range("G133").Value = "FK Gabala"
$ws.Range("H133").Value = 1
$ws.Range("I133").Value = 1
$ws.Range("J133").Value = "D"
$ws.Range("K133").Value = 1.85
$ws.Range("L133").Value = 3.2
$ws.Range("M133").Value = 3.8
$ws.Range("N133").Value = 1.909
$ws.Range("O133").Value = 3.1
$ws.Range("P133").Value = 3.6
$ws.Range("Q133").Value = -0.5
$ws.Range("R133").Value = 1.975
$ws.Range("S133").Value = 1.825
$ws.Range("T133").Value = 2.25
$ws.Range("U133").Value = 1.975
$ws.Range("V133").Value = 1.725
$ws.Range("W133").Value = -1
$ws.Range("X133").Value = 2.1
$ws.Range("Y133").Value = -1
$ws.Range("Z133").Value = -1
$ws.Range("AA133").Value = 0.825
$ws.Range("AB133").Value = -0.5
$ws.Range("AC133").Value = 0.3625

# ---- Row 134 (future fixture: no result yet, missing last two odds cols) ----
$ws.Range("A134").Value = 132
$ws.Range("B134").Value = 7011623
$ws.Range("C134").Value = "Azerbaijan Premier League"
$ws.Range("D134").Value = "Azerbaijan Premier League"
$ws.Range("E134").Value = 45368.5
$ws.Range("F134").Value = "FK Sumqayit"
$ws.Range("G134").Value = "Sabah"
$ws.Range("K134").Value = 2.9
$ws.Range("L134").Value = 3.75
$ws.Range("M134").Value = 2
$ws.Range("N134").Value = 2.75
$ws.Range("O134").Value = 3.5
$ws.Range("P134").Value = 2.15
$ws.Range("Q134").Value = 0.25
$ws.Range("R134").Value = 1.825
$ws.Range("S134").Value = 1.975
$ws.Range("T134").Value = 2.25
$ws.Range("U134").Value = 1.975
$ws.Range("V134").Value = 1.825
$ws.Range("W134").Value = 0
$ws.Range("X134").Value = 0
$ws.Range("Y134").Value = 0
$ws.Range("Z134").Value = 0
$ws.Range("AA134").Value = 0

# Clear the cells that must stay completely empty (no value AND no style)
# for row 134: match result columns (H/I/J) and the last two odds columns
# (AB/AC) were not present in the source diff.
$ws.Range("H134:J134").ClearContents()
$ws.Range("AB134:AC134").ClearContents()
